$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear previously-empty placeholder cells (D/E/F/G) that should no longer be
# present as explicit (empty) cells in rows 2-20.
$ws.Range("D2:G5").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("D14:G15").ClearContents()
$ws.Range("E16:G17").ClearContents()
$ws.Range("D18:G18").ClearContents()
$ws.Range("D19").ClearContents()
$ws.Range("D20:G20").ClearContents()

# Append the new "Address" model rows (21-26).
$ws.Range("A21").Value = "locations"
$ws.Range("B21").Value = "Address"
$ws.Range("C21").Value = "City"
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = "locations"
$ws.Range("F21").Value = "City"
$ws.Range("G21").Value = "GeographicalUnitName"
$ws.Range("H21").Value = $false

$ws.Range("A22").Value = "locations"
$ws.Range("B22").Value = "Address"
$ws.Range("C22").Value = "Country"
$ws.Range("D22").Value = ""
$ws.Range("E22").Value = "locations"
$ws.Range("F22").Value = "Country"
$ws.Range("G22").Value = "GeographicalUnitName"
$ws.Range("H22").Value = $false

$ws.Range("A23").Value = "locations"
$ws.Range("B23").Value = "Address"
$ws.Range("C23").Value = "created_by"
$ws.Range("D23").Value = ""
$ws.Range("E23").Value = "auth"
$ws.Range("F23").Value = "user"
$ws.Range("G23").Value = "username"
$ws.Range("H23").Value = $false

$ws.Range("A24").Value = "locations"
$ws.Range("B24").Value = "Address"
$ws.Range("C24").Value = "updated_by"
$ws.Range("D24").Value = ""
$ws.Range("E24").Value = "auth"
$ws.Range("F24").Value = "user"
$ws.Range("G24").Value = "username"
$ws.Range("H24").Value = $false

$ws.Range("A25").Value = "locations"
$ws.Range("B25").Value = "Address"
$ws.Range("C25").Value = "updated"
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = ""
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = $true

$ws.Range("A26").Value = "locations"
$ws.Range("B26").Value = "Address"
$ws.Range("C26").Value = "created"
$ws.Range("D26").Value = ""
$ws.Range("E26").Value = ""
$ws.Range("F26").Value = ""
$ws.Range("G26").Value = ""
$ws.Range("H26").Value = $true
